# Auto-generated Excel COM-interop edit script
# Updates market-data derived columns (currentAveragePrice*, LevePrice*, LeveProfit*)
# across multiple sheets per the scheduled market-data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 913.7778
$ws.Range("I33").Value = 946.6667
$ws.Range("K33").Value = 946.6667
$ws.Range("M33").Value = -717.6667

$ws.Range("H76").Value = 8996.5
$ws.Range("I76").Value = 8996.5
$ws.Range("K76").Value = 8996.5
$ws.Range("M76").Value = -8681.5

$ws.Range("H79").Value = 8996.5
$ws.Range("I79").Value = 8996.5
$ws.Range("K79").Value = 8996.5
$ws.Range("M79").Value = -7904.5

$ws.Range("H86").Value = 181822580
$ws.Range("I86").Value = 250003950
$ws.Range("K86").Value = 250003950
$ws.Range("M86").Value = -250002827

$ws.Range("H89").Value = 181822580
$ws.Range("I89").Value = 250003950
$ws.Range("K89").Value = 1250019750
$ws.Range("M89").Value = -1250014134

$ws.Range("H125").Value = 1506.1333
$ws.Range("J125").Value = 3820.25
$ws.Range("L125").Value = 34382.25
$ws.Range("N125").Value = -39302.25

$ws.Range("H134").Value = 142434
$ws.Range("J134").Value = 142434
$ws.Range("L134").Value = 142434
$ws.Range("N134").Value = -152574

$ws.Range("H137").Value = 3485.4443
$ws.Range("I137").Value = 4963.815
$ws.Range("K137").Value = 14891.445
$ws.Range("M137").Value = -12341.445

$ws.Range("H141").Value = 16496.889
$ws.Range("I141").Value = 20518.889
$ws.Range("J141").Value = 8452.888999999999
$ws.Range("K141").Value = 61556.667
$ws.Range("L141").Value = 25358.667
$ws.Range("M141").Value = -56376.667
$ws.Range("N141").Value = -35718.667

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4766.3423
$ws.Range("I32").Value = 4896.448
$ws.Range("K32").Value = 4896.448
$ws.Range("M32").Value = -4609.448

$ws.Range("H61").Value = 2868.46
$ws.Range("I61").Value = 2573.9512
$ws.Range("K61").Value = 2573.9512
$ws.Range("M61").Value = -2361.9512

$ws.Range("H74").Value = 1757.5358
$ws.Range("I74").Value = 1168.44
$ws.Range("J74").Value = 6666.6665
$ws.Range("K74").Value = 1168.44
$ws.Range("L74").Value = 6666.6665
$ws.Range("M74").Value = -294.4400000000001
$ws.Range("N74").Value = -8414.666499999999

$ws.Range("H77").Value = 1757.5358
$ws.Range("I77").Value = 1168.44
$ws.Range("J77").Value = 6666.6665
$ws.Range("K77").Value = 5842.200000000001
$ws.Range("L77").Value = 33333.3325
$ws.Range("M77").Value = -1474.200000000001
$ws.Range("N77").Value = -42069.3325

$ws.Range("H110").Value = 1434
$ws.Range("I110").Value = 1170.0625
$ws.Range("K110").Value = 1170.0625
$ws.Range("M110").Value = 874.9375

$ws.Range("H114").Value = 87599.39999999999
$ws.Range("J114").Value = 87599.39999999999
$ws.Range("L114").Value = 87599.39999999999
$ws.Range("N114").Value = -96277.39999999999

$ws.Range("H136").Value = 2868.46
$ws.Range("I136").Value = 2573.9512
$ws.Range("K136").Value = 7721.8536
$ws.Range("M136").Value = -5171.8536

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()

$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()

$ws.Range("H134").Value = 3449.95
$ws.Range("I134").Value = 2986.8655
$ws.Range("J134").Value = 6460
$ws.Range("K134").Value = 8960.5965
$ws.Range("L134").Value = 19380
$ws.Range("M134").Value = -6425.5965
$ws.Range("N134").Value = -24450

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H53").Value = 40099.832
$ws.Range("J53").Value = 40099.832
$ws.Range("L53").Value = 40099.832
$ws.Range("N53").Value = -41313.832

$ws.Range("H141").Value = 92234.266
$ws.Range("J141").Value = 98222.69
$ws.Range("L141").Value = 98222.69
$ws.Range("N141").Value = -108582.69

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 78.88
$ws.Range("I2").Value = 41.333332
$ws.Range("K2").Value = 247.999992
$ws.Range("M2").Value = -134.999992

$ws.Range("H38").Value = 324.22726
$ws.Range("I38").Value = 54.454544
$ws.Range("J38").Value = 594
$ws.Range("K38").Value = 163.363632
$ws.Range("L38").Value = 1782
$ws.Range("M38").Value = 183.636368
$ws.Range("N38").Value = -2476

$ws.Range("H97").Value = 304.7857
$ws.Range("I97").Value = 166.42857
$ws.Range("J97").Value = 443.14285
$ws.Range("K97").Value = 499.28571
$ws.Range("L97").Value = 1329.42855
$ws.Range("M97").Value = -3.285709999999995
$ws.Range("N97").Value = -2321.42855

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 25000
$ws.Range("I43").Value = 20000
$ws.Range("K43").Value = 20000
$ws.Range("M43").Value = -19849

$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()

$ws.Range("H122").Value = 10222.556
$ws.Range("I122").Value = 9403.4
$ws.Range("J122").Value = 11246.5
$ws.Range("K122").Value = 28210.2
$ws.Range("L122").Value = 33739.5
$ws.Range("M122").Value = -25760.2
$ws.Range("N122").Value = -38639.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H5").Value = 1000000
$ws.Range("I5").Value = 1000000
$ws.Range("K5").Value = 1000000
$ws.Range("M5").Value = -999887

$ws.Range("H7").Value = 8250
$ws.Range("I7").Value = 8250
$ws.Range("K7").Value = 8250
$ws.Range("M7").Value = -8138

$ws.Range("H100").Value = 76926320
$ws.Range("I100").Value = 166670190
$ws.Range("K100").Value = 166670190
$ws.Range("M100").Value = -166669649

$ws.Range("H126").Value = 8250
$ws.Range("I126").Value = 8250
$ws.Range("K126").Value = 24750
$ws.Range("M126").Value = -22280

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 7005712
$ws.Range("I81").Value = 7589355
$ws.Range("K81").Value = 15178710
$ws.Range("M81").Value = -15177649

$ws.Range("H84").Value = 7005712
$ws.Range("I84").Value = 7589355
$ws.Range("K84").Value = 75893550
$ws.Range("M84").Value = -75888246

$ws.Range("H132").Value = 2374.32
$ws.Range("I132").Value = 2409.311
$ws.Range("J132").Value = 2059.4
$ws.Range("K132").Value = 7227.933000000001
$ws.Range("L132").Value = 6178.200000000001
$ws.Range("M132").Value = -4697.933000000001
$ws.Range("N132").Value = -11238.2

$ws.Range("H135").Value = 250713.67
$ws.Range("J135").Value = 250713.67
$ws.Range("L135").Value = 250713.67
$ws.Range("N135").Value = -260853.67
